# Update the "K" (strikeouts) column (column G) values for rows 2-16.
# These values were regenerated to reflect K (strikeouts) instead of the
# previous "Strike#" (total strikes thrown) metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 5
    4  = 3
    5  = 2
    6  = 5
    7  = 2
    8  = 4
    9  = 4
    10 = 4
    11 = 3
    12 = 0
    13 = 2
    14 = 8
    15 = 4
    16 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
